$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1583.3334
$ws.Range("J29").Value = 4150
$ws.Range("L29").Value = 12450
$ws.Range("N29").Value = -13012

$ws.Range("H32").Value = 1000
$ws.Range("I32").Value = 1000
$ws.Range("K32").Value = 1000
$ws.Range("M32").Value = -674

$ws.Range("H40").Value = 2319.5
$ws.Range("J40").Value = 1999
$ws.Range("L40").Value = 1999
$ws.Range("N40").Value = -2349

$ws.Range("H41").Value = 425
$ws.Range("J41").Value = 300
$ws.Range("L41").Value = 300
$ws.Range("N41").Value = -1180

$ws.Range("H82").Value = 2845
$ws.Range("I82").Value = 2845
$ws.Range("K82").Value = 8535
$ws.Range("M82").Value = -8129

$ws.Range("H85").Value = 2845
$ws.Range("I85").Value = 2845
$ws.Range("K85").Value = 8535
$ws.Range("M85").Value = -7131

$ws.Range("H99").Value = 47.5
$ws.Range("I99").Value = 47.5
$ws.Range("K99").Value = 142.5
$ws.Range("M99").Value = 1355.5

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H37").Value = 32630.6
$ws.Range("I37").Value = 3157
$ws.Range("K37").Value = 3157
$ws.Range("M37").Value = -2884

$ws.Range("H62").Value = 44166.668
$ws.Range("J62").Value = 44166.668
$ws.Range("L62").Value = 44166.668
$ws.Range("N62").Value = -45414.668

$ws.Range("H65").Value = 44166.668
$ws.Range("J65").Value = 44166.668
$ws.Range("L65").Value = 132500.004
$ws.Range("N65").Value = -138740.004

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H92").Value = 45066.168
$ws.Range("J92").Value = 45066.168
$ws.Range("L92").Value = 45066.168
$ws.Range("N92").Value = -50058.168

$ws.Range("H94").Value = 30000
$ws.Range("J94").Value = 30000
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31802

$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 11380
$ws.Range("J75").Value = 20000
$ws.Range("L75").Value = 20000
$ws.Range("N75").Value = -21872

$ws.Range("H78").Value = 11380
$ws.Range("J78").Value = 20000
$ws.Range("L78").Value = 60000
$ws.Range("N78").Value = -69360

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws.Range("H94").Value = 775
$ws.Range("I94").Value = 775
$ws.Range("K94").Value = 775
$ws.Range("M94").Value = -324

$ws.Range("H134").Value = 7555.3335
$ws.Range("I134").Value = 4499.75
$ws.Range("K134").Value = 13499.25
$ws.Range("M134").Value = -10964.25

$ws.Range("H135").Value = 93374
$ws.Range("J135").Value = 93374
$ws.Range("L135").Value = 93374
$ws.Range("N135").Value = -103514

$ws.Range("H140").Value = 161308.53
$ws.Range("I140").Value = 308750
$ws.Range("K140").Value = 308750
$ws.Range("M140").Value = -303570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 270.83334
$ws.Range("I22").Value = 370
$ws.Range("K22").Value = 370
$ws.Range("M22").Value = -20

$ws.Range("H132").Value = 2917.25
$ws.Range("I132").Value = 2335.3333
$ws.Range("K132").Value = 7005.999899999999
$ws.Range("M132").Value = -4475.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 833
$ws.Range("I5").Value = 750
$ws.Range("K5").Value = 2250
$ws.Range("M5").Value = -2138

$ws.Range("H23").Value = 389.5
$ws.Range("I23").Value = 371.5
$ws.Range("K23").Value = 1114.5
$ws.Range("M23").Value = -879.5

$ws.Range("H34").Value = 921.8570999999999
$ws.Range("I34").Value = 200
$ws.Range("J34").Value = 1210.6
$ws.Range("K34").Value = 600
$ws.Range("L34").Value = 3631.8
$ws.Range("M34").Value = -516
$ws.Range("N34").Value = -3799.8

$ws.Range("H55").Value = 298.33334
$ws.Range("J55").Value = 347.5
$ws.Range("L55").Value = 1042.5
$ws.Range("N55").Value = -1396.5

$ws.Range("H135").Value = 833
$ws.Range("I135").Value = 750
$ws.Range("K135").Value = 6750
$ws.Range("M135").Value = -4215

$ws.Range("H139").Value = 2423.1667
$ws.Range("I139").Value = 1769.5
$ws.Range("J139").Value = 2750
$ws.Range("K139").Value = 5308.5
$ws.Range("L139").Value = 8250
$ws.Range("M139").Value = -168.5
$ws.Range("N139").Value = -18530

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 6500
$ws.Range("I33").Value = 6500
$ws.Range("K33").Value = 6500
$ws.Range("M33").Value = -6248

$ws.Range("H36").Value = 117
$ws.Range("I36").Value = 117
$ws.Range("K36").Value = 117
$ws.Range("M36").Value = 368

$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()

$ws.Range("H122").Value = 1709.4
$ws.Range("I122").Value = 886.75
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 2660.25
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -210.25
$ws.Range("N122").Value = -19900

$ws.Range("H132").Value = 2263.8
$ws.Range("I132").Value = 2127.6
$ws.Range("K132").Value = 6382.799999999999
$ws.Range("M132").Value = -3852.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 60000
$ws.Range("I68").Value = 60000
$ws.Range("K68").Value = 60000
$ws.Range("M68").Value = -59251

$ws.Range("H71").Value = 60000
$ws.Range("I71").Value = 60000
$ws.Range("K71").Value = 300000
$ws.Range("M71").Value = -296256

$ws.Range("H76").Value = 22500
$ws.Range("J76").Value = 22500
$ws.Range("L76").Value = 22500
$ws.Range("N76").Value = -23176

$ws.Range("H79").Value = 22500
$ws.Range("J79").Value = 22500
$ws.Range("L79").Value = 22500
$ws.Range("N79").Value = -24840

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 52094
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

$ws.Range("H94").Value = 25000
$ws.Range("J94").Value = 25000
$ws.Range("L94").Value = 25000
$ws.Range("N94").Value = -26802

$ws.Range("H103").Value = 19602
$ws.Range("J103").Value = 19602
$ws.Range("L103").Value = 19602
$ws.Range("N103").Value = -21946

$ws.Range("H135").Value = 80000
$ws.Range("J135").Value = 80000
$ws.Range("L135").Value = 80000
$ws.Range("N135").Value = -90140
